$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 31.14997866666667
$ws.Range("H2").Value = 93.44993600000001
$ws.Range("I2").Value = 0.4621739036316256
$ws.Range("J2").Value = 0.4621739036316256
$ws.Range("M2").Value = 0.4200680000000001
$ws.Range("N2").Value = 1.260204
$ws.Range("O2").Value = 0.01461602726853518
$ws.Range("P2").Value = 0.01461602726853518
$ws.Range("Q2").Value = 13.08510923854934
$ws.Range("R2").Value = 117.765983146944
$ws.Range("S2").Value = 0.00675514637828519
$ws.Range("T2").Value = 0.00675514637828519

$ws.Range("G3").Value = 31.14997866666667
$ws.Range("H3").Value = 93.44993600000001
$ws.Range("I3").Value = 0.4621739036316256
$ws.Range("J3").Value = 0.4621739036316256
$ws.Range("O3").Value = 0.3267187263271011
$ws.Range("P3").Value = 0.3267187263271011
$ws.Range("Q3").Value = 292.4974170972712
$ws.Range("R3").Value = 2632.47675387544
$ws.Range("S3").Value = 0.1510008691361491
$ws.Range("T3").Value = 0.1510008691361491

$ws.Range("G4").Value = 31.14997866666667
$ws.Range("H4").Value = 93.44993600000001
$ws.Range("I4").Value = 0.4621739036316256
$ws.Range("J4").Value = 0.4621739036316256
$ws.Range("M4").Value = 18.93019133333333
$ws.Range("N4").Value = 56.79057399999999
$ws.Range("O4").Value = 0.6586652464043636
$ws.Range("P4").Value = 0.6586652464043636
$ws.Range("Q4").Value = 589.6750561892514
$ws.Range("R4").Value = 5307.075505703264
$ws.Range("S4").Value = 0.3044178881171913
$ws.Range("T4").Value = 0.3044178881171913

$ws.Range("G5").Value = 18.94069966666667
$ws.Range("H5").Value = 56.822099
$ws.Range("I5").Value = 0.2810241764892454
$ws.Range("J5").Value = 0.2810241764892454
$ws.Range("M5").Value = 0.4200680000000001
$ws.Range("N5").Value = 1.260204
$ws.Range("O5").Value = 0.01461602726853518
$ws.Range("P5").Value = 0.01461602726853518
$ws.Range("Q5").Value = 7.956381827577334
$ws.Range("R5").Value = 71.60743644819601
$ws.Range("S5").Value = 0.004107457026684453
$ws.Range("T5").Value = 0.004107457026684453

$ws.Range("G6").Value = 18.94069966666667
$ws.Range("H6").Value = 56.822099
$ws.Range("I6").Value = 0.2810241764892454
$ws.Range("J6").Value = 0.2810241764892454
$ws.Range("O6").Value = 0.3267187263271011
$ws.Range("P6").Value = 0.3267187263271011
$ws.Range("S6").Value = 0.09181586100968873
$ws.Range("T6").Value = 0.09181586100968873

$ws.Range("G7").Value = 18.94069966666667
$ws.Range("H7").Value = 56.822099
$ws.Range("I7").Value = 0.2810241764892454
$ws.Range("J7").Value = 0.2810241764892454
$ws.Range("M7").Value = 18.93019133333333
$ws.Range("N7").Value = 56.79057399999999
$ws.Range("O7").Value = 0.6586652464043636
$ws.Range("P7").Value = 0.6586652464043636
$ws.Range("Q7").Value = 358.5510686772028
$ws.Range("R7").Value = 3226.959618094826
$ws.Range("S7").Value = 0.1851008584528722
$ws.Range("T7").Value = 0.1851008584528722

$ws.Range("G8").Value = 14.86848
$ws.Range("H8").Value = 44.60544
$ws.Range("I8").Value = 0.2206044349565553
$ws.Range("J8").Value = 0.2206044349565553
$ws.Range("M8").Value = 0.4200680000000001
$ws.Range("N8").Value = 1.260204
$ws.Range("O8").Value = 0.01461602726853518
$ws.Range("P8").Value = 0.01461602726853518
$ws.Range("Q8").Value = 6.245772656640001
$ws.Range("R8").Value = 56.21195390976001
$ws.Range("S8").Value = 0.003224360436884807
$ws.Range("T8").Value = 0.003224360436884807

$ws.Range("G9").Value = 14.86848
$ws.Range("H9").Value = 44.60544
$ws.Range("I9").Value = 0.2206044349565553
$ws.Range("J9").Value = 0.2206044349565553
$ws.Range("O9").Value = 0.3267187263271011
$ws.Range("P9").Value = 0.3267187263271011
$ws.Range("Q9").Value = 139.6146059264
$ws.Range("R9").Value = 1256.5314533376
$ws.Range("S9").Value = 0.07207560001111557
$ws.Range("T9").Value = 0.07207560001111557

$ws.Range("G10").Value = 14.86848
$ws.Range("H10").Value = 44.60544
$ws.Range("I10").Value = 0.2206044349565553
$ws.Range("J10").Value = 0.2206044349565553
$ws.Range("M10").Value = 18.93019133333333
$ws.Range("N10").Value = 56.79057399999999
$ws.Range("O10").Value = 0.6586652464043636
$ws.Range("P10").Value = 0.6586652464043636
$ws.Range("Q10").Value = 281.4631712358399
$ws.Range("R10").Value = 2533.16854112256
$ws.Range("S10").Value = 0.1453044745085549
$ws.Range("T10").Value = 0.1453044745085549

$ws.Range("G11").Value = 2.439668
$ws.Range("H11").Value = 7.319004
$ws.Range("I11").Value = 0.03619748492257375
$ws.Range("J11").Value = 0.03619748492257375
$ws.Range("M11").Value = 0.4200680000000001
$ws.Range("N11").Value = 1.260204
$ws.Range("O11").Value = 0.01461602726853518
$ws.Range("P11").Value = 0.01461602726853518
$ws.Range("Q11").Value = 1.024826457424
$ws.Range("R11").Value = 9.223438116816
$ws.Range("S11").Value = 0.0005290634266807289
$ws.Range("T11").Value = 0.0005290634266807289

$ws.Range("G12").Value = 2.439668
$ws.Range("H12").Value = 7.319004
$ws.Range("I12").Value = 0.03619748492257375
$ws.Range("J12").Value = 0.03619748492257375
$ws.Range("O12").Value = 0.3267187263271011
$ws.Range("P12").Value = 0.3267187263271011
$ws.Range("Q12").Value = 22.90841339607333
$ws.Range("R12").Value = 206.17572056466
$ws.Range("S12").Value = 0.01182639617014774
$ws.Range("T12").Value = 0.01182639617014774

$ws.Range("G13").Value = 2.439668
$ws.Range("H13").Value = 7.319004
$ws.Range("I13").Value = 0.03619748492257375
$ws.Range("J13").Value = 0.03619748492257375
$ws.Range("M13").Value = 18.93019133333333
$ws.Range("N13").Value = 56.79057399999999
$ws.Range("O13").Value = 0.6586652464043636
$ws.Range("P13").Value = 0.6586652464043636
$ws.Range("Q13").Value = 46.18338202981065
$ws.Range("R13").Value = 415.6504382682959
$ws.Range("S13").Value = 0.02384202532574528
$ws.Range("T13").Value = 0.02384202532574528
